# io-showcase.xlsx: add new Web ("rightClick(locator)") and Xml
# ("insertAfter(xml,xpath,content,var)", "insertBefore(xml,xpath,content,var)",
# "replaceIn(xml,xpath,content,var)") command-reference entries to the hidden
# '#system' lookup sheet, keeping each column alphabetically sorted, and grow
# the "web" / "xml" named ranges to cover the newly-added rows.

function Insert-SortedValue($ws, $col, $lastRow, $insertRow, $value) {
    # Push every value from $insertRow..$lastRow down by one row, then drop
    # the new value into the now-empty $insertRow. Walking bottom-up means
    # each cell is read before it gets overwritten.
    for ($r = $lastRow; $r -ge $insertRow; $r--) {
        $ws.Cells.Item($r + 1, $col).Value2 = $ws.Cells.Item($r, $col).Value2
    }
    $ws.Cells.Item($insertRow, $col).Value2 = $value
    return $lastRow + 1
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- column AC ("xml" system lookup list, was $AC$2:$AC$18) -----------------
# Insert the 3 new xml commands in their correct alphabetical slots.
$lastRowAC = 18
$lastRowAC = Insert-SortedValue $ws 29 $lastRowAC 13 "insertAfter(xml,xpath,content,var)"
$lastRowAC = Insert-SortedValue $ws 29 $lastRowAC 14 "insertBefore(xml,xpath,content,var)"
$lastRowAC = Insert-SortedValue $ws 29 $lastRowAC 18 "replaceIn(xml,xpath,content,var)"

# --- column X ("web" system lookup list, was $X$2:$X$122) -------------------
# Insert the new rightClick(locator) command in its alphabetical slot.
$lastRowX = 122
$lastRowX = Insert-SortedValue $ws 24 $lastRowX 79 "rightClick(locator)"

# --- keep the named ranges in sync with the new list lengths -----------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$X`$2:`$X`$$lastRowX"
    }
    if ($n.Name -eq "xml") {
        $n.RefersTo = "='#system'!`$AC`$2:`$AC`$$lastRowAC"
    }
}
